$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 48097424
$ws.Range("I40").Value = 1554.3636
$ws.Range("J40").Value = 101002880
$ws.Range("K40").Value = 1554.3636
$ws.Range("L40").Value = 101002880
$ws.Range("M40").Value = -1379.3636
$ws.Range("N40").Value = -101003230
# Row 88
$ws.Range("H88").Value = 6761974.5
$ws.Range("I88").Value = 1093.6666
$ws.Range("J88").Value = 10142415
$ws.Range("K88").Value = 1093.6666
$ws.Range("L88").Value = 10142415
$ws.Range("M88").Value = -687.6666
$ws.Range("N88").Value = -10143227
# Row 91
$ws.Range("H91").Value = 6761974.5
$ws.Range("I91").Value = 1093.6666
$ws.Range("J91").Value = 10142415
$ws.Range("K91").Value = 1093.6666
$ws.Range("L91").Value = 10142415
$ws.Range("M91").Value = 310.3334
$ws.Range("N91").Value = -10145223
# Row 97
$ws.Range("H97").Value = 143572430
$ws.Range("J97").Value = 143572430
$ws.Range("L97").Value = 430717290
$ws.Range("N97").Value = -430718282
# Row 137
$ws.Range("H137").Value = 1219.3518
$ws.Range("I137").Value = 995.125
$ws.Range("J137").Value = 1860
$ws.Range("K137").Value = 2985.375
$ws.Range("L137").Value = 5580
$ws.Range("M137").Value = -435.375
$ws.Range("N137").Value = -10680

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 19667742
$ws.Range("I45").Value = 25718918
$ws.Range("J45").Value = 1424.75
$ws.Range("K45").Value = 25718918
$ws.Range("L45").Value = 1424.75
$ws.Range("M45").Value = -25718541
$ws.Range("N45").Value = -2178.75
# Row 74
$ws.Range("H74").Value = 1539.25
$ws.Range("I74").Value = 1620.5
$ws.Range("J74").Value = 970.5
$ws.Range("K74").Value = 1620.5
$ws.Range("L74").Value = 970.5
$ws.Range("M74").Value = -746.5
$ws.Range("N74").Value = -2718.5
# Row 77
$ws.Range("H77").Value = 1539.25
$ws.Range("I77").Value = 1620.5
$ws.Range("J77").Value = 970.5
$ws.Range("K77").Value = 8102.5
$ws.Range("L77").Value = 4852.5
$ws.Range("M77").Value = -3734.5
$ws.Range("N77").Value = -13588.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1448.766
$ws.Range("I31").Value = 1127.9706
$ws.Range("J31").Value = 2287.7693
$ws.Range("K31").Value = 1127.9706
$ws.Range("L31").Value = 2287.7693
$ws.Range("M31").Value = -832.9706000000001
$ws.Range("N31").Value = -2877.7693
# Row 34
$ws.Range("H34").Value = 1448.766
$ws.Range("I34").Value = 1127.9706
$ws.Range("J34").Value = 2287.7693
$ws.Range("K34").Value = 1127.9706
$ws.Range("L34").Value = 2287.7693
$ws.Range("M34").Value = -925.9706000000001
$ws.Range("N34").Value = -2691.7693
# Row 62
$ws.Range("H62").Value = 4022.8
$ws.Range("I62").Value = 2773
$ws.Range("J62").Value = 4647.7
$ws.Range("K62").Value = 2773
$ws.Range("L62").Value = 4647.7
$ws.Range("M62").Value = -2149
$ws.Range("N62").Value = -5895.7
# Row 65
$ws.Range("H65").Value = 4022.8
$ws.Range("I65").Value = 2773
$ws.Range("J65").Value = 4647.7
$ws.Range("K65").Value = 13865
$ws.Range("L65").Value = 23238.5
$ws.Range("M65").Value = -10745
$ws.Range("N65").Value = -29478.5
# Row 132
$ws.Range("H132").Value = 22223902
$ws.Range("I132").Value = 1894.4445
$ws.Range("J132").Value = 55556916
$ws.Range("K132").Value = 5683.333500000001
$ws.Range("L132").Value = 166670748
$ws.Range("M132").Value = -3153.333500000001
$ws.Range("N132").Value = -166675808

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1356.5
$ws.Range("I46").Value = 668.8889
$ws.Range("J46").Value = 2594.2
$ws.Range("K46").Value = 668.8889
$ws.Range("L46").Value = 2594.2
$ws.Range("M46").Value = -480.8889
$ws.Range("N46").Value = -2970.2
# Row 68
$ws.Range("H68").Value = 1525.375
$ws.Range("I68").Value = 1540.6
$ws.Range("J68").Value = 1500
$ws.Range("K68").Value = 1540.6
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -791.5999999999999
$ws.Range("N68").Value = -2998
# Row 69
$ws.Range("H69").Value = 33382
$ws.Range("J69").Value = 33382
$ws.Range("L69").Value = 33382
$ws.Range("N69").Value = -35004
# Row 71
$ws.Range("H71").Value = 1525.375
$ws.Range("I71").Value = 1540.6
$ws.Range("J71").Value = 1500
$ws.Range("K71").Value = 7703
$ws.Range("L71").Value = 7500
$ws.Range("M71").Value = -3959
$ws.Range("N71").Value = -14988
# Row 72
$ws.Range("H72").Value = 33382
$ws.Range("J72").Value = 33382
$ws.Range("L72").Value = 100146
$ws.Range("N72").Value = -108258
# Row 136
$ws.Range("H136").Value = 40819064
$ws.Range("I136").Value = 6806022.5
$ws.Range("J136").Value = 142858190
$ws.Range("K136").Value = 20418067.5
$ws.Range("L136").Value = 428574570
$ws.Range("M136").Value = -20415517.5
$ws.Range("N136").Value = -428579670

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5389.15
$ws.Range("I62").Value = 5225.8184
$ws.Range("J62").Value = 5588.778
$ws.Range("K62").Value = 5225.8184
$ws.Range("L62").Value = 5588.778
$ws.Range("M62").Value = -4601.8184
$ws.Range("N62").Value = -6836.778
# Row 65
$ws.Range("H65").Value = 5389.15
$ws.Range("I65").Value = 5225.8184
$ws.Range("J65").Value = 5588.778
$ws.Range("K65").Value = 26129.092
$ws.Range("L65").Value = 27943.89
$ws.Range("M65").Value = -23009.092
$ws.Range("N65").Value = -34183.89
# Row 81
$ws.Range("H81").Value = 898.4286
$ws.Range("I81").Value = 813.53845
$ws.Range("K81").Value = 1627.0769
$ws.Range("M81").Value = -566.0769
# Row 84
$ws.Range("H84").Value = 898.4286
$ws.Range("I84").Value = 813.53845
$ws.Range("K84").Value = 8135.3845
$ws.Range("M84").Value = -2831.3845
# Row 132
$ws.Range("H132").Value = 2454.36
$ws.Range("I132").Value = 1249.7397
$ws.Range("J132").Value = 5711.2964
$ws.Range("K132").Value = 3749.2191
$ws.Range("L132").Value = 17133.8892
$ws.Range("M132").Value = -1219.2191
$ws.Range("N132").Value = -22193.8892
# Row 133
$ws.Range("H133").Value = 42707.5
$ws.Range("J133").Value = 42707.5
$ws.Range("L133").Value = 42707.5
$ws.Range("N133").Value = -52827.5
# Row 136
$ws.Range("H136").Value = 9263742
$ws.Range("I136").Value = 13519051
$ws.Range("K136").Value = 40557153
$ws.Range("M136").Value = -40554603
